$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-09-01 Sunday" "2024-09-02 Monday"

Replace-Text "94×52=" "39×74="
Replace-Text "29×76=" "30×78="
Replace-Text "37×79=" "60×99="
Replace-Text "36×40=" "49×74="
Replace-Text "60×28=" "49×85="
Replace-Text "45×63=" "75×26="
Replace-Text "49×29=" "12×67="
Replace-Text "37×70=" "70×26="
Replace-Text "69×87=" "23×86="
Replace-Text "52×47=" "57×85="
Replace-Text "60×84=" "64×48="
Replace-Text "28×71=" "43×60="
Replace-Text "97×27=" "59×54="
Replace-Text "67×94=" "18×60="
Replace-Text "19×17=" "70×85="
Replace-Text "24×20=" "60×55="
Replace-Text "11×38=" "76×28="
Replace-Text "77×17=" "49×76="
Replace-Text "78×38=" "13×16="
Replace-Text "26×38=" "68×24="
Replace-Text "75×82=" "56×76="
Replace-Text "87×44=" "18×17="
Replace-Text "99×31=" "58×73="
Replace-Text "12×64=" "26×32="
Replace-Text "30×86=" "72×33="
